$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 35, shifting rows 35-45 down to 36-46
$ws.Rows("35").Insert()

# Populate new row 35 with data (copy of common fields, new specific values)
$ws.Cells.Item(35, 1).Value = 5
$ws.Cells.Item(35, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(35, 3).Value = "Maule"
$ws.Cells.Item(35, 4).Value = 44785
$ws.Cells.Item(35, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 5).Value = 7
$ws.Cells.Item(35, 6).Value = 100112040
$ws.Cells.Item(35, 7).Value = "Cilantro"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 150
$ws.Cells.Item(35, 11).Value = 14000
$ws.Cells.Item(35, 12).Value = 14000
$ws.Cells.Item(35, 13).Value = 14000
$ws.Cells.Item(35, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(35, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(35, 16).Value = 389
$ws.Cells.Item(35, 17).Value = 36
$ws.Cells.Item(35, 18).Value = "Hortaliza"
